$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$summary = $wb.Worksheets.Item("Summary")

# Rows where only the R column value needs to be cleared
$rOnlyRows = @(4, 5, 14, 15)
foreach ($r in $rOnlyRows) {
    $ws.Range("R$r").Value = $null
}

# Rows where P, R, S, T columns need to be cleared (Q stays)
$fullRows = @(8, 10, 12, 18, 24, 26, 28)
foreach ($r in $fullRows) {
    $ws.Range("P$r").Value = $null
    $ws.Range("R$r").Value = $null
    $ws.Range("S$r").Value = $null
    $ws.Range("T$r").Value = $null
}

# Update the Summary sheet's "Total Backward GEMM (G_ops)" value
$summary.Range("B5").Value = 0.004194304

# Update the conditional formatting threshold on R1:R29 (Backward GEMM ops)
# from 16777216 down to 2097152
$fcs = $ws.Range("R1:R29").FormatConditions
$fc = $fcs.Item(1)
$fc.Formula1 = "=2097152"
